$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new commit "pull added & pick repaired", 1 hour
$ws.Range("C16").Value = "pull added & pick repaired"
$ws.Range("G16").Value = 1

# Row 17: new commit "npc deleted & merchant & talk instruction added", 1.5 hours
$ws.Range("C17").Value = "npc deleted & merchant & talk instruction added"
$ws.Range("G17").Value = 1.5

# Match the author's formatting: the commit-name cells use the same
# style as the rest of the filled commit rows (C5:C15)
$ws.Range("C16").Style = $ws.Range("C15").Style
$ws.Range("C17").Style = $ws.Range("C15").Style

# Move the active selection, as left by the author after editing
$ws.Range("E18").Select()
